{"js": "// \"made multiple bug fixes, including banking logic\"\n//  1) Strike through the three existing \"bug\" paragraphs (they're being\n//     marked as fixed / crossed out), applied at both the paragraph-mark\n//     level and on every run in the paragraph.\n//  2) Append a new paragraph reporting another bug (\"Neil Armstrong\" is\n//     misspelled as \"amrstrong\"), left un-struck, with the misspelled word\n//     flagged the same way Word's proofer flags the other typos in the doc.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Strike through every one of the original paragraphs.\nfor (const paragraph of paragraphs.items) {\n  paragraph.font.strikeThrough = true;\n}\nawait context.sync();\n\n// 2) Add the new \"Neil Armstrong,amrstrong\" paragraph after the last one,\n//    rebuilding the exact run / proofErr split (two runs, \"Neil \" and\n//    \"Armstrong,amrstrong\", the second wrapped in spellStart/spellEnd\n//    proof-error markers) via a literal OOXML insert so the markup matches\n//    the rest of the document's spell-checked runs.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst newParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nconst newRange = newParagraph.getRange();\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Neil </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>Armstrong,amrstrong</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nnewRange.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# \"made multiple bug fixes, including banking logic\"\n#  1) Strike through the three existing \"bug\" paragraphs (they're being\n#     marked as fixed / crossed out), applied at both the paragraph-mark\n#     level and on every run in the paragraph.\n#  2) Append a new paragraph reporting another bug (\"Neil Armstrong\" is\n#     misspelled as \"amrstrong\"), left un-struck, with the misspelled word\n#     flagged the same way Word's proofer flags the other typos in the doc.\n\n$d = $word.ActiveDocument\n\n# 1) Strike through every one of the original paragraphs.\n$originalCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $originalCount; $i++) {\n    $d.Paragraphs.Item($i).Range.Font.StrikeThrough = 1\n}\n\n# 2) Add a new empty paragraph after the last existing one, then fill it in\n#    via a literal OOXML insert so we can rebuild the exact run / proofErr\n#    split (two runs, \"Neil \" and \"Armstrong,amrstrong\", the second wrapped\n#    in spellStart/spellEnd proof-error markers) that matches the rest of\n#    the document's spell-checked runs.\n$lastPara = $d.Paragraphs.Last\n$tail = $lastPara.Range\n$tail.Collapse(0)\n$tail.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newRange = $newPara.Range\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + `\n    '<w:p>' + `\n    '<w:r><w:t xml:space=\"preserve\">Neil </w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellStart\"/>' + `\n    '<w:r><w:t>Armstrong,amrstrong</w:t></w:r>' + `\n    '<w:proofErr w:type=\"spellEnd\"/>' + `\n    '</w:p>' + `\n    '</w:body>' + `\n    '</w:document>' + `\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n# \"Replace\" tells Word to replace just the (empty) new paragraph's own\n# content/mark with the supplied XML instead of inserting alongside it \u2014\n# without it an extra blank paragraph is left behind after the insert.\n$newRange.InsertXML($flatOpc, \"Replace\")\n"}
